# Update the "Förändrad" date column (C2:C8) from 2023-10-13 (45212) to 2023-10-22 (45221)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C8").Value = 45221
